$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 92, pushing the existing rows 92-119 down to 93-120.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with this week's data. The columns that
# are constant for every data row in this sheet (A, B, C, E, F, G, H, I, N, Q,
# R) are copied from the surrounding rows; the rest are the new figures.
$ws.Cells.Item(92, 1).Value = 10
$ws.Cells.Item(92, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(92, 3).Value = "La Araucanía"
$ws.Cells.Item(92, 4).Value = 44798
$ws.Cells.Item(92, 5).Value = 9
$ws.Cells.Item(92, 6).Value = 100112035
$ws.Cells.Item(92, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 90
$ws.Cells.Item(92, 11).Value = 24000
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 13).Value = 24556
$ws.Cells.Item(92, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 2456
$ws.Cells.Item(92, 17).Value = 10
$ws.Cells.Item(92, 18).Value = "Hortaliza"
